$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "temporalidad" column (B): was a dimension, is now re-curated as a measure.
$ws.Range("B2").Value = "iaest-measure:temporalidad"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").Clear()

# "sector-descripcion" column (H): was a dimension, is now re-curated as a measure.
$ws.Range("H2").Value = "iaest-measure:sector-descripcion"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("H5").Clear()
